$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (sample-size / id header values)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON - meanEMG values)
$ws.Range("B2").Value = 15.238634862404243
$ws.Range("C2").Value = 10.711959559265367
$ws.Range("D2").Value = 15.411627666745767
$ws.Range("E2").Value = 8.2976607373479681

# Row 3 (STR - meanEMG values)
$ws.Range("B3").Value = 13.034818003848427
$ws.Range("C3").Value = 13.314304289333535
$ws.Range("D3").Value = 11.881139235523573
$ws.Range("E3").Value = 14.297854565260399

# Update selection to match the new active range used in the authored workbook
$ws.Range("B1:E3").Select()
